$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarter labels for rows 6-11 (column A), continuing Q0..Q9
$labels = @("Q4", "Q5", "Q6", "Q7", "Q8", "Q9")

# New data rows (B..G) matching rows 6-11
$data = @(
  @(0.2046491862038776, 0.7536088930238598, 3.651661730967382, 1.910932162837651, 1.922972664121156, 42),
  @(0.2392227887939334, 0.7909065342463171, 3.777102677640282, 1.943476955778041, 1.95265767950657, 41),
  @(0.2155957076965697, 0.7869099249525766, 3.856332724670826, 1.963754751660916, 1.976749783912624, 40),
  @(0.239545290171405, 0.8122047914116975, 3.974152750272352, 1.993527714949645, 2.004954792540174, 39),
  @(0.2295130447590046, 0.80090761193645, 4.057843671499801, 2.01440901296132, 2.028155615889948, 38),
  @(0.189254522310732, 0.7743474122443053, 3.997164790338612, 1.999291071939904, 2.017767317389937, 37)
)

# Updated values for existing rows 2-5 (B..G)
$updated = @{
  2 = @(0.1121240238528402, 0.7618801882577865, 3.665753451068926, 1.914615745017503, 1.932450120671928, 46)
  3 = @(0.2384546409039757, 0.7169237981258901, 3.429598611108179, 1.851917549759756, 1.857253653004148, 45)
  4 = @(0.1554256882402988, 0.746261807642007, 3.542524738013336, 1.882159594193153, 1.897416690034527, 44)
  5 = @(0.2282236696078792, 0.7489212578072509, 3.600547596828088, 1.897510895048586, 1.906029573040004, 43)
}

foreach ($r in 2..5) {
  $vals = $updated[$r]
  $ws.Cells.Item($r, 2).Value = $vals[0]
  $ws.Cells.Item($r, 3).Value = $vals[1]
  $ws.Cells.Item($r, 4).Value = $vals[2]
  $ws.Cells.Item($r, 5).Value = $vals[3]
  $ws.Cells.Item($r, 6).Value = $vals[4]
  $ws.Cells.Item($r, 7).Value = $vals[5]
}

# Reference style cell to copy formatting from (A2, which already carries the
# bordered/bold/centered style used for the quarter-label column)
$styleSource = $ws.Range("A2")
$styleSource.Copy()

for ($i = 0; $i -lt $labels.Length; $i++) {
  $r = 6 + $i
  $cellA = $ws.Cells.Item($r, 1)
  $cellA.Value = $labels[$i]
  $cellA.PasteSpecial(-4122)  # xlPasteFormats

  $vals = $data[$i]
  $ws.Cells.Item($r, 2).Value = $vals[0]
  $ws.Cells.Item($r, 3).Value = $vals[1]
  $ws.Cells.Item($r, 4).Value = $vals[2]
  $ws.Cells.Item($r, 5).Value = $vals[3]
  $ws.Cells.Item($r, 6).Value = $vals[4]
  $ws.Cells.Item($r, 7).Value = $vals[5]
}
